$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - LinearRegression (name unchanged, B2 unchanged)
$ws.Range("C2").Value = 3650464472355304
$ws.Range("D2").Value = 3650464472355304

# Row 3 - RandomForestRegressor (name unchanged)
$ws.Range("B3").Value = 61662351217919.29
$ws.Range("C3").Value = 57880570928473.56
$ws.Range("D3").Value = 657521942064222.4

# Row 4 - GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.03872732500726751
$ws.Range("C4").Value = 0.03780677747620992
$ws.Range("D4").Value = 184481545137872

# Row 5 - AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 130251668181522.4
$ws.Range("C5").Value = 21514723638725.91
$ws.Range("D5").Value = 247892735450383.5
